$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current row 628 (2026/12/29 火 13 9),
# shifting all rows 628-669 down to 629-670.
$ws.Rows.Item(628).Insert()

# Fill the new row 628 with the inserted record: 2026/01/15 木 9 30.
# Column A holds a date-like string ("2026/01/15") that must be stored as
# literal text (matching every other date cell in the sheet), not coerced
# into a date serial number, so force a text number format before writing
# the value, then restore the cell to the workbook's normal (unstyled) look.
$cellA = $ws.Cells.Item(628, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/01/15"
$cellA.Style = "Normal"

$ws.Cells.Item(628, 2).Value = "木"
$ws.Cells.Item(628, 3).Value = 9
$ws.Cells.Item(628, 4).Value = 30
